# Fruta / hortaliza, semanal
# Insert a new weekly data row at row 186 (pushing the existing rows 186-203
# down to 187-204) on the single data sheet, then populate the new row with
# the latest week's record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 186..203 down to 187..204, inserting a fresh blank row 186.
$ws.Rows.Item(186).Insert()

# Populate the newly inserted row 186 with the new weekly record.
$ws.Cells.Item(186, 1).Value2  = 10
$ws.Cells.Item(186, 2).Value2  = "Vega Modelo de Temuco"
$ws.Cells.Item(186, 3).Value2  = "La Araucanía"
$ws.Cells.Item(186, 4).Value2  = 45021
$ws.Cells.Item(186, 5).Value2  = 9
$ws.Cells.Item(186, 6).Value2  = 100112031
$ws.Cells.Item(186, 7).Value2  = "Poroto verde"
$ws.Cells.Item(186, 8).Value2  = "Brío"
$ws.Cells.Item(186, 9).Value2  = "Primera"
$ws.Cells.Item(186, 10).Value2 = 200
$ws.Cells.Item(186, 11).Value2 = 1000
$ws.Cells.Item(186, 12).Value2 = 1000
$ws.Cells.Item(186, 13).Value2 = 1000
$ws.Cells.Item(186, 14).Value2 = "$/kilo"
$ws.Cells.Item(186, 15).Value2 = "Región de La Araucanía"
$ws.Cells.Item(186, 16).Value2 = 1000
$ws.Cells.Item(186, 17).Value2 = 1
$ws.Cells.Item(186, 18).Value2 = "Hortaliza"
